$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) and Volume(1h) (E) columns keep their literal text formatting
# (no auto-conversion of numeric-looking strings to numbers, preserving trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "44.114.91"
$ws.Range("E2").Value = "  +0.84%  "

$ws.Range("D3").Value = "2.246.44"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "272.62"
$ws.Range("E5").Value = "  +5.14%  "

$ws.Range("D6").Value = "87.87"
$ws.Range("E6").Value = "  +11.00%  "

$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  -0.08%  "

$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  +0.69%  "

$ws.Range("D10").Value = "45.39"
$ws.Range("E10").Value = "  +4.80%  "

$ws.Range("D11").Value = "0.0929"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").Value = "7.68"
$ws.Range("E12").Value = "  +8.16%  "

$ws.Range("E13").Value = "  +1.54%  "

$ws.Range("D14").Value = "2.588.32"
$ws.Range("E14").Value = "  +0.22%  "

$ws.Range("D15").Value = "14.97"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").Value = "2.259.64"
$ws.Range("E16").Value = "  +0.66%  "

$ws.Range("D17").Value = "0.798"
$ws.Range("E17").Value = "  +0.24%  "

$ws.Range("D18").Value = "44.053.63"
$ws.Range("E18").Value = "  +0.79%  "

$ws.Range("E19").Value = "  -0.49%  "

$ws.Range("D20").Value = "6.00"
$ws.Range("E20").Value = "  -1.27%  "

$ws.Range("D21").Value = "70.56"
$ws.Range("E21").Value = "  -1.46%  "

$ws.Range("D22").Value = "2.37"
$ws.Range("E22").Value = "  +2.01%  "

$ws.Range("D23").Value = "234.45"
$ws.Range("E23").Value = "  +0.52%  "

$ws.Range("D24").Value = "8.76"
$ws.Range("E24").Value = "  -7.12%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("D26").Value = "2.54"
$ws.Range("E26").Value = "  +13.56%  "

$ws.Range("D27").Value = "10.87"
$ws.Range("E27").Value = "  +0.26%  "

$ws.Range("D28").Value = "3.55"
$ws.Range("E28").Value = "  +5.92%  "

$ws.Range("D29").Value = "2.31"
$ws.Range("E29").Value = "  +5.14%  "

$ws.Range("D30").Value = "39.78"
$ws.Range("E30").Value = "  -5.65%  "

$ws.Range("D31").Value = "174.53"
$ws.Range("E31").Value = "  +0.78%  "

$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "0.0910"
$ws.Range("E32").Value = "  +3.75%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "20.89"
$ws.Range("E33").Value = "  +1.39%  "

$ws.Range("E34").Value = "  +2.19%  "

$ws.Range("E35").Value = "  +0.77%  "

$ws.Range("D37").Value = "0.0355"
$ws.Range("E37").Value = "  -3.11%  "

$ws.Range("D38").Value = "4.36"
$ws.Range("E38").Value = "  -2.95%  "

$ws.Range("E39").Value = "  +21.20%  "

$ws.Range("D40").Value = "2.20"
$ws.Range("E40").Value = "  +2.97%  "

$ws.Range("D41").Value = "12.55"
$ws.Range("E41").Value = "  -4.86%  "

$ws.Range("D42").Value = "64.56"
$ws.Range("E42").Value = "  +4.31%  "

$ws.Range("E43").Value = "  +0.83%  "

$ws.Range("D44").Value = "5.46"
$ws.Range("E44").Value = "  +1.15%  "

$ws.Range("D45").Value = "8.51"
$ws.Range("E45").Value = "  -0.50%  "

$ws.Range("D46").Value = "0.0989"
$ws.Range("E46").Value = "  +0.17%  "

$ws.Range("D47").Value = "100.39"
$ws.Range("E47").Value = "  -3.93%  "

$ws.Range("D48").Value = "1.20"
$ws.Range("E48").Value = "  +4.03%  "

$ws.Range("D49").Value = "1.14"
$ws.Range("E49").Value = "  +0.98%  "

$ws.Range("D50").Value = "0.429"
$ws.Range("E50").Value = "  -8.83%  "

$ws.Range("D51").Value = "1.48"
$ws.Range("E51").Value = "  -0.49%  "
